$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 614, shifting the existing rows 614:671 down to 615:672
$ws.Rows("614:614").Insert()

# Populate the newly inserted row 614 with the new weekly observation
$ws.Range("A614").Value = 6
$ws.Range("B614").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C614").Value = "Metropolitana"
$ws.Range("D614").Value = 45127
$ws.Range("E614").Value = 13
$ws.Range("F614").Value = "Fruta"
$ws.Range("G614").Value = 100107
$ws.Range("H614").Value = "Otros"
$ws.Range("I614").Value = 100107011
$ws.Range("J614").Value = "Tuna"
$ws.Range("K614").Value = "Sin especificar"
$ws.Range("L614").Value = "Extra (doble especial)"
$ws.Range("M614").Value = 120
$ws.Range("N614").Value = 29000
$ws.Range("O614").Value = 30000
$ws.Range("P614").Value = 29500
$ws.Range("Q614").Value = "`$/caja 18 kilos"
$ws.Range("R614").Value = "Provincia de Melipilla"
$ws.Range("S614").Value = 1639
$ws.Range("T614").Value = 18
